$wb = $excel.ActiveWorkbook

# Sheets: 1=TestSuite, 2=LoginPageTest, 3=AdminPageTest
$wsLogin = $wb.Worksheets.Item(2)
$wsAdmin = $wb.Worksheets.Item(3)

# --- LoginPageTest (sheet2): selection becomes a block A1:B2 ---
$wsLogin.Activate()
$wsLogin.Range("A1:B2").Select()

# --- AdminPageTest (sheet3): insert the UserName/Password columns taken
# from LoginPageTest in front of the existing User/SearchCriteria/RunMode
# columns, shifting the old data from A:C to C:E. ---
$wsAdmin.Activate()
$wsAdmin.Columns("A:B").Insert()

$wsAdmin.Range("A1").Value = "UserName"
$wsAdmin.Range("B1").Value = "Password"
$wsAdmin.Range("A2").Value = "Admin"
$wsAdmin.Range("B2").Value = "admin123"

# Copy the header/data formatting (fill + borders) from the existing
# columns so the new ones look consistent with the rest of the table.
$wsAdmin.Range("C1:C2").Copy()
$wsAdmin.Range("A1:B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize the new columns to fit their contents (closest achievable fit
# given this host's column-width rounding granularity).
$wsAdmin.Columns("A").ColumnWidth = 9.5
$wsAdmin.Columns("B").ColumnWidth = 8.666666666666666

$wsAdmin.Range("D1").Select()
